$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4612464904785156
$ws.Range("B1").Value = 0.3770990669727325
$ws.Range("C1").Value = 3.732804775238037
$ws.Range("D1").Value = 3.466596841812134
$ws.Range("E1").Value = 0.9372816681861877
